$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) date serial from 46074 to 46075
# for every data row (rows 2 through 39).
$ws.Range("C2:C39").Value2 = 46075
